$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B3 and B4 from "No" to "Yes"
$ws.Range("B3").Value = "Yes"
$ws.Range("B4").Value = "Yes"

# Update the selection on Sheet1 to I13
$ws.Activate()
$ws.Range("I13").Select()
